$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bordered style used by the existing data rows (row 2) down into
# the new rows 4-7 so the new data matches the look of the existing rows.
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A4:C7").PasteSpecial(-4122) | Out-Null

# Row 4: adminMentor / abc123 / Invalid
$ws.Range("A4").Value = "adminMentor"
$ws.Range("B4").Value = "abc123"
$ws.Range("C4").Value = "Invalid"

# Row 5: studentAdmin / xyz123 / Invalid
$ws.Range("A5").Value = "studentAdmin"
$ws.Range("B5").Value = "xyz123"
$ws.Range("C5").Value = "Invalid"

# Row 6: sangeeta / sangeeta123 / Valid
$ws.Range("A6").Value = "sangeeta"
$ws.Range("B6").Value = "sangeeta123"
$ws.Range("C6").Value = "Valid"

# Row 7: sweetapal / sweeta123 / Valid
$ws.Range("A7").Value = "sweetapal"
$ws.Range("B7").Value = "sweeta123"
$ws.Range("C7").Value = "Valid"

# Match the saved selection state from the target workbook.
$ws.Range("C6").Select() | Out-Null
